$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1719038817005545
$ws.Range("C2").Value = 0.6118299445471349
$ws.Range("J2").Value = 0.011090573012939
$ws.Range("P2").Value = 0.121996303142329
$ws.Range("S2").Value = 0.08317929759704251
$ws.Range("B3").Value = 0.002849002849002849
$ws.Range("C3").Value = 0.04273504273504274
$ws.Range("J3").Value = 0.02849002849002849
$ws.Range("P3").Value = 0.7264957264957265
$ws.Range("S3").Value = 0.1994301994301994
$ws.Range("J4").Value = 0.05617977528089887
$ws.Range("O4").Value = 0.01123595505617977
$ws.Range("P4").Value = 0.7191011235955056
$ws.Range("S4").Value = 0.2134831460674157
$ws.Range("B6").Value = 0.04375
$ws.Range("D6").Value = 0.02083333333333333
$ws.Range("E6").Value = 0.00625
$ws.Range("F6").Value = 0.08125
$ws.Range("J6").Value = 0.2458333333333333
$ws.Range("O6").Value = 0.01041666666666667
$ws.Range("Q6").Value = 0.1729166666666667
$ws.Range("R6").Value = 0.06666666666666667
$ws.Range("S6").Value = 0.3520833333333334
$ws.Range("B7").Value = 0.08724832214765101
$ws.Range("D7").Value = 0.01118568232662192
$ws.Range("E7").Value = 0.002237136465324385
$ws.Range("F7").Value = 0.06935123042505593
$ws.Range("J7").Value = 0.1252796420581656
$ws.Range("O7").Value = 0.01118568232662192
$ws.Range("Q7").Value = 0.1722595078299776
$ws.Range("R7").Value = 0.08724832214765101
$ws.Range("S7").Value = 0.4340044742729307
$ws.Range("B8").Value = 0.09895833333333333
$ws.Range("D8").Value = 0.01979166666666667
$ws.Range("E8").Value = 0.002083333333333333
$ws.Range("F8").Value = 0.06041666666666667
$ws.Range("J8").Value = 0.109375
$ws.Range("O8").Value = 0.01979166666666667
$ws.Range("Q8").Value = 0.178125
$ws.Range("R8").Value = 0.08333333333333333
$ws.Range("S8").Value = 0.428125
$ws.Range("B9").Value = 0.07555555555555556
$ws.Range("D9").Value = 0.02
$ws.Range("F9").Value = 0.07111111111111111
$ws.Range("J9").Value = 0.1066666666666667
$ws.Range("O9").Value = 0.01333333333333333
$ws.Range("Q9").Value = 0.1844444444444444
$ws.Range("R9").Value = 0.08222222222222222
$ws.Range("S9").Value = 0.4466666666666667
$ws.Range("B10").Value = 0.1018518518518518
$ws.Range("D10").Value = 0.01892109500805153
$ws.Range("E10").Value = 0.002012882447665056
$ws.Range("F10").Value = 0.07971014492753623
$ws.Range("J10").Value = 0.1123188405797101
$ws.Range("O10").Value = 0.01288244766505636
$ws.Range("Q10").Value = 0.2121578099838969
$ws.Range("R10").Value = 0.08695652173913043
$ws.Range("S10").Value = 0.3731884057971014
$ws.Range("G11").Value = 0.1430769230769231
$ws.Range("J11").Value = 0.09230769230769231
$ws.Range("K11").Value = 0.2030769230769231
$ws.Range("L11").Value = 0.5507692307692308
$ws.Range("S11").Value = 0.01076923076923077
$ws.Range("G12").Value = 0.7890410958904109
$ws.Range("J12").Value = 0.1561643835616438
$ws.Range("K12").Value = 0.005479452054794521
$ws.Range("L12").Value = 0.01917808219178082
$ws.Range("S12").Value = 0.03013698630136986
$ws.Range("G13").Value = 0.6885245901639344
$ws.Range("J13").Value = 0.2950819672131147
$ws.Range("S13").Value = 0.01639344262295082
$ws.Range("F15").Value = 0.02237136465324385
$ws.Range("H15").Value = 0.203579418344519
$ws.Range("I15").Value = 0.0738255033557047
$ws.Range("J15").Value = 0.3400447427293065
$ws.Range("K15").Value = 0.06263982102908278
$ws.Range("M15").Value = 0.02237136465324385
$ws.Range("N15").Value = 0.002237136465324385
$ws.Range("O15").Value = 0.06711409395973154
$ws.Range("S15").Value = 0.2058165548098434
$ws.Range("F16").Value = 0.01308900523560209
$ws.Range("H16").Value = 0.2094240837696335
$ws.Range("I16").Value = 0.0968586387434555
$ws.Range("J16").Value = 0.3821989528795812
$ws.Range("K16").Value = 0.09162303664921466
$ws.Range("M16").Value = 0.02356020942408377
$ws.Range("N16").Value = 0.002617801047120419
$ws.Range("O16").Value = 0.05235602094240838
$ws.Range("S16").Value = 0.1282722513089005
$ws.Range("F17").Value = 0.02330508474576271
$ws.Range("H17").Value = 0.173728813559322
$ws.Range("I17").Value = 0.1038135593220339
$ws.Range("J17").Value = 0.4025423728813559
$ws.Range("K17").Value = 0.1165254237288136
$ws.Range("M17").Value = 0.02754237288135593
$ws.Range("O17").Value = 0.06991525423728813
$ws.Range("S17").Value = 0.0826271186440678
$ws.Range("F18").Value = 0.0175
$ws.Range("H18").Value = 0.175
$ws.Range("I18").Value = 0.1175
$ws.Range("J18").Value = 0.3875
$ws.Range("K18").Value = 0.1
$ws.Range("M18").Value = 0.0275
$ws.Range("N18").Value = 0.0025
$ws.Range("O18").Value = 0.0675
$ws.Range("S18").Value = 0.105
$ws.Range("F19").Value = 0.01349206349206349
$ws.Range("H19").Value = 0.2222222222222222
$ws.Range("I19").Value = 0.09523809523809523
$ws.Range("J19").Value = 0.3583333333333333
$ws.Range("K19").Value = 0.119047619047619
$ws.Range("M19").Value = 0.02738095238095238
$ws.Range("N19").Value = 0.001190476190476191
$ws.Range("O19").Value = 0.06984126984126984
$ws.Range("S19").Value = 0.09325396825396826
